# Applies the edits described by the commit:
#  "importPage can upload multiple files. Now remain validate"
#
# Concretely (derived from the OOXML diff):
#  - A5 changes from the plain number 2590081110 to the equivalent text string
#  - C5 gets a new value "Dja " (trailing space)
#  - D5 gets a new value "dmak"
#  - E5 (previously "asd") is cleared
#  - D10 gets a new value "susu"
#  - D14 gets a new value "pass it"
#  - F12 gets a new value "get loose"
#  - the sheet view scrolls so column E is left-most and F12 becomes the
#    active/selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- A5: was a bare number 2590081110, becomes the text "2590081110" -------
$A5 = $ws.Range("A5")
$A5.NumberFormat = "@"
$A5.Value = "2590081110"

# --- new / changed cell values ---------------------------------------------
$ws.Range("C5").Value = "Dja "
$ws.Range("D5").Value = "dmak"
$ws.Range("E5").ClearContents()
$ws.Range("D10").Value = "susu"
$ws.Range("D14").Value = "pass it"
$ws.Range("F12").Value = "get loose"

# --- sheet view: scroll so E1 is the top-left cell, select F12 -------------
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F12").Select() | Out-Null
